# Generate Report for Handback
#
# Once a handback package has been produced, record the handback
# target file, the generated handback (.xlf) file and the handback
# timestamp for every localized language sheet, and flip the overall
# status from "Ready for handoff" to "Handed back: in sync with en-US"
# everywhere that status is shown.

$wb = $excel.ActiveWorkbook

$statusText  = "Handed back: in sync with en-US"
$targetFile  = "4066dd2f-145f-4bc4-9a29-b8822e9b16ac.md"
$targetUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6f24fd9ef516f1944f4a07fbf850af948bcc2b5a/e2e/4066dd2f-145f-4bc4-9a29-b8822e9b16ac.md"

# ---- Overview sheet: status shown per-language ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusText
$ov.Range("F2").Value = $statusText
$ov.Range("E3").Value = $statusText
$ov.Range("F3").Value = $statusText

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

$zh.Range("J2").Value = $targetFile
$zh.Hyperlinks.Add($zh.Range("J2"), $targetUrl, "", "", $targetFile)
$zh.Range("K2").Value = "4066dd2f-145f-4bc4-9a29-b8822e9b16ac.1f5c1ef0c073e683166b6b5c438544bd79b4898e.zh-cn.xlf"
$zh.Range("L2").Value = "2017-02-09 16:05:41"

$zh.Range("J3").Value = $targetFile
$zh.Hyperlinks.Add($zh.Range("J3"), $targetUrl, "", "", $targetFile)
$zh.Range("K3").Value = "4066dd2f-145f-4bc4-9a29-b8822e9b16ac.1f5c1ef0c073e683166b6b5c438544bd79b4898e.zh-cn.xlf"
$zh.Range("L3").Value = "2017-02-09 16:05:41"

$zh.Columns.Item(3).ColumnWidth = 29.9777050018311
$zh.Columns.Item(10).ColumnWidth = 40
$zh.Columns.Item(11).ColumnWidth = 40

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

$de.Range("J2").Value = $targetFile
$de.Hyperlinks.Add($de.Range("J2"), $targetUrl, "", "", $targetFile)
$de.Range("K2").Value = "4066dd2f-145f-4bc4-9a29-b8822e9b16ac.1f5c1ef0c073e683166b6b5c438544bd79b4898e.de-de.xlf"
$de.Range("L2").Value = "2017-02-09 16:06:07"

$de.Range("J3").Value = $targetFile
$de.Hyperlinks.Add($de.Range("J3"), $targetUrl, "", "", $targetFile)
$de.Range("K3").Value = "4066dd2f-145f-4bc4-9a29-b8822e9b16ac.1f5c1ef0c073e683166b6b5c438544bd79b4898e.de-de.xlf"
$de.Range("L3").Value = "2017-02-09 16:06:07"

$de.Columns.Item(3).ColumnWidth = 29.9777050018311
$de.Columns.Item(10).ColumnWidth = 40
$de.Columns.Item(11).ColumnWidth = 40

# ---- Overview sheet: widen the per-language status columns too ----
$ov.Columns.Item(5).ColumnWidth = 29.9777050018311
$ov.Columns.Item(6).ColumnWidth = 29.9777050018311

Write-Output "done"
